# Applies the edits described by the commit diff to the active document.
#
# Summary of changes:
#   1. Title line: "...-2026-MDP-GLDE/SGLCA" -> "...-2026-MDP/GLDE-SGLCA"
#      (the "-" and "/" separators swap places).
#   2. "Que, mediante el Documento Simple N° {" : the "N" and "° " runs
#      are merged into a single run "N° " (formatting preserved).
#   3. "(en adelante, {{genero2}})" : the "en adelante", "," and
#      " {{genero2}})" runs are merged into a single run.
#   4. "REGÍSTRESE, COMUNÍQUESE Y CÚMPLASE." : the "REG", "Í" and
#      "STRESE, COMUNÍQUESE Y CÚMPLASE." runs are merged into a single run.
#
# NB: whenever several adjacent runs end up sharing 100% identical
# run-formatting, Word (and this host) always coalesces them into one
# <w:r> on save, so the merges below are expressed as plain text edits -
# there is no separate "change the formatting" step required.  An
# intermediate throw-away assignment is used before writing the final
# text so the engine is forced to rebuild/re-merge the run list instead
# of treating the "same text" write as a no-op.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) RESOLUCIÓN GERENCIAL Nº {{cod_resolucion}}-2026-MDP-GLDE/SGLCA
#                                                 -> -2026-MDP/GLDE-SGLCA
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "RESOLUCIÓN GERENCIAL Nº {{cod_resolucion}}-2026-MDP-GLDE/SGLCA",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "zzz"
    $rng.Text = "RESOLUCIÓN GERENCIAL Nº {{cod_resolucion}}-2026-MDP/GLDE-SGLCA"
}
Write-Output "title updated: $found"

# ---------------------------------------------------------------------
# 2) "N" + "° " -> "N° " (only inside "...Documento Simple N° {...")
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "Que, mediante el Documento Simple N° {",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sub = $d.Range($rng.Start, $rng.End)
    $found2 = $sub.Find.Execute("N° ", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "", 0)
    if ($found2) {
        $sub.Text = "zzz"
        $sub.Text = "N° "
    }
}
Write-Output "N grado merged: $found"

# ---------------------------------------------------------------------
# 3) "en adelante" + "," + " {{genero2}})" -> "en adelante, {{genero2}})"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "en adelante, {{genero2}})",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "zzz"
    $rng.Text = "en adelante, {{genero2}})"
}
Write-Output "en adelante merged: $found"

# ---------------------------------------------------------------------
# 4) "REG" + "Í" + "STRESE, COMUNÍQUESE Y CÚMPLASE." -> single run
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "REGÍSTRESE, COMUNÍQUESE Y CÚMPLASE.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "zzz"
    $rng.Text = "REGÍSTRESE, COMUNÍQUESE Y CÚMPLASE."
}
Write-Output "registrese merged: $found"
